$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sensor Map")
$ws.Range("A1").Value = "test"
